# Regenerate merged AHB files
#
# 1) Rename header row columns: "_old" -> "_FV2310", "_new" -> "_FV2404"
# 2) Convert the data range A1:U58 into an Excel Table (ListObject) named "Table1"
#    while preserving the existing header-row formatting (bold/fill/border)
#    instead of letting Excel capture it as a table "header override" dxf.
# 3) Freeze the header row (pane split below row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header cells -----------------------------------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2) Turn the used range into a Table ---------------------------------
# Stash the header row's existing formatting outside the used range so we
# can put it back after ListObjects.Add(); otherwise Excel would record the
# pre-existing bold/fill header look as a per-table "headerRowDxfId" override.
$headerRange = $ws.Range("A1:U1")
$stashRange = $ws.Range("A60:U60")

$headerRange.Copy()
$stashRange.PasteSpecial(-4122) # xlPasteFormats
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U58")
$list = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$list.Name = "Table1"
$list.TableStyle = ""

$stashRange.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats
$stashRange.Clear()
$excel.CutCopyMode = $false

# --- 3) Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
